$d = $word.ActiveDocument

# Find the paragraph that contains the "(give the name of the log file)" text
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*(give the name of the log file)*") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after it
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()

# Copy the formatting (run + paragraph mark formatting) from the source
# paragraph onto the freshly inserted one, so fonts/bold/color/size match.
$srcRange = $target.Range.Duplicate()
$srcRange.MoveEnd(1, -1) | Out-Null
$newPara.Range.FormattedText = $srcRange.FormattedText

# Replace the copied text with the new log-file name, keeping the
# paragraph mark (and its formatting) intact.
$newRange = $newPara.Range.Duplicate()
$newRange.MoveEnd(1, -1) | Out-Null
$newRange.Text = "projectlog.xlsx / log.pdf"
